$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "clientes": add two new client rows (4 and 5)
# ---------------------------------------------------------------------------
$clientes = $wb.Worksheets.Item("clientes")

# Row 4
$clientes.Cells.Item(4,1).Value = "activo"
$clientes.Cells.Item(4,2).Value = "16.742.249-7"
$clientes.Cells.Item(4,3).Value = "Isaias Beroiza Mora"
$clientes.Cells.Item(4,4).Value = "colaco sn km3 parcela 9"
$clientes.Cells.Item(4,5).Value = "Calbuco"
$clientes.Cells.Item(4,6).NumberFormat = "@"
$clientes.Cells.Item(4,6).Value = "88809703"
$clientes.Cells.Item(4,7).Value = "por buscar"
$clientes.Cells.Item(4,8).Value = "ok"

# Row 5
$clientes.Cells.Item(5,1).Value = "activo"
$clientes.Cells.Item(5,2).Value = "17673326-8"
$clientes.Cells.Item(5,3).Value = "Maria Jose Rodriguez"
$clientes.Cells.Item(5,4).Value = "colaco sn km3 parcela 9"
$clientes.Cells.Item(5,5).Value = "ca"
$clientes.Cells.Item(5,6).NumberFormat = "@"
$clientes.Cells.Item(5,6).Value = "88809704"
$clientes.Cells.Item(5,7).Value = "por buscar"
$clientes.Cells.Item(5,8).Value = "ok"

# ---------------------------------------------------------------------------
# Sheet "ruta_actual": add the same two clients to the current route (rows 4 and 5)
# ---------------------------------------------------------------------------
$ruta = $wb.Worksheets.Item("ruta_actual")

# Row 4
$ruta.Cells.Item(4,1).NumberFormat = "@"
$ruta.Cells.Item(4,1).Value = "20240707"
$ruta.Cells.Item(4,2).Value = 2
$ruta.Cells.Item(4,3).Value = "16.742.249-7"
$ruta.Cells.Item(4,4).Value = "Isaias Beroiza Mora"
$ruta.Cells.Item(4,5).Value = "colaco sn km3 parcela 9"
$ruta.Cells.Item(4,6).Value = "Calbuco"
$ruta.Cells.Item(4,7).NumberFormat = "@"
$ruta.Cells.Item(4,7).Value = "88809703"
$ruta.Cells.Item(4,8).Value = "por buscar"
$ruta.Cells.Item(4,9).Value = "ok"

# Row 5
$ruta.Cells.Item(5,1).NumberFormat = "@"
$ruta.Cells.Item(5,1).Value = "20240707"
$ruta.Cells.Item(5,2).Value = 3
$ruta.Cells.Item(5,3).Value = "17673326-8"
$ruta.Cells.Item(5,4).Value = "Maria Jose Rodriguez"
$ruta.Cells.Item(5,5).Value = "colaco sn km3 parcela 9"
$ruta.Cells.Item(5,6).Value = "ca"
$ruta.Cells.Item(5,7).NumberFormat = "@"
$ruta.Cells.Item(5,7).Value = "88809704"
$ruta.Cells.Item(5,8).Value = "por buscar"
$ruta.Cells.Item(5,9).Value = "ok"
